# Auto-generated data-driven update of computed profit columns (H:N)
# across multiple sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "ALC" = @(
        @(@{ Cell = "H42"; Value = 629.8889 }; @{ Cell = "I42"; Value = 435.8 }; @{ Cell = "J42"; Value = 872.5 }; @{ Cell = "K42"; Value = 1307.4 }; @{ Cell = "L42"; Value = 2617.5 }; @{ Cell = "M42"; Value = -1077.4 }; @{ Cell = "N42"; Value = -3077.5 })
        @(@{ Cell = "H88"; Value = 2283 }; @{ Cell = "I88"; Value = 0 }; @{ Cell = "J88"; Value = 2283 }; @{ Cell = "K88"; Value = 0 }; @{ Cell = "L88"; Value = 2283 }; @{ Cell = "M88"; Value = $null }; @{ Cell = "N88"; Value = -3095 })
        @(@{ Cell = "H91"; Value = 2283 }; @{ Cell = "I91"; Value = 0 }; @{ Cell = "J91"; Value = 2283 }; @{ Cell = "K91"; Value = 0 }; @{ Cell = "L91"; Value = 2283 }; @{ Cell = "M91"; Value = $null }; @{ Cell = "N91"; Value = -5091 })
        @(@{ Cell = "H129"; Value = 2783 }; @{ Cell = "I129"; Value = 2574 }; @{ Cell = "K129"; Value = 7722 }; @{ Cell = "M129"; Value = -2722 })
        @(@{ Cell = "H137"; Value = 2115.7036 }; @{ Cell = "J137"; Value = 4172.385 }; @{ Cell = "L137"; Value = 12517.155 }; @{ Cell = "N137"; Value = -17617.155 })
        @(@{ Cell = "H138"; Value = 4980.4287 }; @{ Cell = "I138"; Value = 4332.3335 }; @{ Cell = "J138"; Value = 5157.1816 }; @{ Cell = "K138"; Value = 12997.0005 }; @{ Cell = "L138"; Value = 15471.5448 }; @{ Cell = "M138"; Value = -7857.000499999998 }; @{ Cell = "N138"; Value = -25751.5448 })
    )
    "ARM" = @(
        @(@{ Cell = "H32"; Value = 19358.191 }; @{ Cell = "I32"; Value = 19304.791 }; @{ Cell = "K32"; Value = 19304.791 }; @{ Cell = "M32"; Value = -19017.791 })
        @(@{ Cell = "H61"; Value = 1834.2727 }; @{ Cell = "I61"; Value = 1663.375 }; @{ Cell = "K61"; Value = 1663.375 }; @{ Cell = "M61"; Value = -1451.375 })
        @(@{ Cell = "H74"; Value = 1820.0741 }; @{ Cell = "I74"; Value = 1243.3684 }; @{ Cell = "J74"; Value = 3189.75 }; @{ Cell = "K74"; Value = 1243.3684 }; @{ Cell = "L74"; Value = 3189.75 }; @{ Cell = "M74"; Value = -369.3684000000001 }; @{ Cell = "N74"; Value = -4937.75 })
        @(@{ Cell = "H77"; Value = 1820.0741 }; @{ Cell = "I77"; Value = 1243.3684 }; @{ Cell = "J77"; Value = 3189.75 }; @{ Cell = "K77"; Value = 6216.842000000001 }; @{ Cell = "L77"; Value = 15948.75 }; @{ Cell = "M77"; Value = -1848.842000000001 }; @{ Cell = "N77"; Value = -24684.75 })
        @(@{ Cell = "H122"; Value = 3720.4 }; @{ Cell = "I122"; Value = 4204 }; @{ Cell = "J122"; Value = 2995 }; @{ Cell = "K122"; Value = 12612 }; @{ Cell = "L122"; Value = 8985 }; @{ Cell = "M122"; Value = -10162 }; @{ Cell = "N122"; Value = -13885 })
        @(@{ Cell = "H136"; Value = 1834.2727 }; @{ Cell = "I136"; Value = 1663.375 }; @{ Cell = "K136"; Value = 4990.125 }; @{ Cell = "M136"; Value = -2440.125 })
    )
    "BSM" = @(
        @(@{ Cell = "H86"; Value = 3253.842 }; @{ Cell = "I86"; Value = 2972 }; @{ Cell = "J86"; Value = 5649.5 }; @{ Cell = "K86"; Value = 2972 }; @{ Cell = "L86"; Value = 5649.5 }; @{ Cell = "M86"; Value = -1849 }; @{ Cell = "N86"; Value = -7895.5 })
        @(@{ Cell = "H89"; Value = 3253.842 }; @{ Cell = "I89"; Value = 2972 }; @{ Cell = "J89"; Value = 5649.5 }; @{ Cell = "K89"; Value = 14860 }; @{ Cell = "L89"; Value = 28247.5 }; @{ Cell = "M89"; Value = -9244 }; @{ Cell = "N89"; Value = -39479.5 })
        @(@{ Cell = "H99"; Value = 1499.6666 }; @{ Cell = "I99"; Value = 1250 }; @{ Cell = "J99"; Value = 1999 }; @{ Cell = "K99"; Value = 1250 }; @{ Cell = "L99"; Value = 1999 }; @{ Cell = "M99"; Value = 248 }; @{ Cell = "N99"; Value = -4995 })
        @(@{ Cell = "H107"; Value = 1028.9048 }; @{ Cell = "I107"; Value = 867.2222 }; @{ Cell = "K107"; Value = 867.2222 }; @{ Cell = "M107"; Value = 1052.7778 })
        @(@{ Cell = "H124"; Value = 66666 }; @{ Cell = "J124"; Value = 66666 }; @{ Cell = "L124"; Value = 66666 }; @{ Cell = "N124"; Value = -76486 })
    )
    "CRP" = @(
        @(@{ Cell = "H31"; Value = 2323.8215 }; @{ Cell = "I31"; Value = 2323.8215 }; @{ Cell = "K31"; Value = 2323.8215 }; @{ Cell = "M31"; Value = -2028.8215 })
        @(@{ Cell = "H34"; Value = 2323.8215 }; @{ Cell = "I34"; Value = 2323.8215 }; @{ Cell = "K34"; Value = 2323.8215 }; @{ Cell = "M34"; Value = -2121.8215 })
        @(@{ Cell = "H99"; Value = 9271.333000000001 }; @{ Cell = "I99"; Value = 9090.5 }; @{ Cell = "K99"; Value = 9090.5 }; @{ Cell = "M99"; Value = -7592.5 })
        @(@{ Cell = "H126"; Value = 9271.333000000001 }; @{ Cell = "I126"; Value = 9090.5 }; @{ Cell = "K126"; Value = 27271.5 }; @{ Cell = "M126"; Value = -24801.5 })
        @(@{ Cell = "H141"; Value = 58784.668 }; @{ Cell = "J141"; Value = 61555 }; @{ Cell = "L141"; Value = 61555 }; @{ Cell = "N141"; Value = -71915 })
    )
    "CUL" = @(
        @(@{ Cell = "H2"; Value = 2057721.2 }; @{ Cell = "J2"; Value = 89.833336 }; @{ Cell = "L2"; Value = 539.000016 }; @{ Cell = "N2"; Value = -765.000016 })
        @(@{ Cell = "H69"; Value = 1612.7142 }; @{ Cell = "I69"; Value = 1263.3334 }; @{ Cell = "J69"; Value = 1874.75 }; @{ Cell = "K69"; Value = 3790.0002 }; @{ Cell = "L69"; Value = 5624.25 }; @{ Cell = "M69"; Value = -2979.0002 }; @{ Cell = "N69"; Value = -7246.25 })
        @(@{ Cell = "H72"; Value = 1612.7142 }; @{ Cell = "I72"; Value = 1263.3334 }; @{ Cell = "J72"; Value = 1874.75 }; @{ Cell = "K72"; Value = 11370.0006 }; @{ Cell = "L72"; Value = 16872.75 }; @{ Cell = "M72"; Value = -7314.000599999999 }; @{ Cell = "N72"; Value = -24984.75 })
        @(@{ Cell = "H131"; Value = 2030.4166 }; @{ Cell = "J131"; Value = 2088.0908 }; @{ Cell = "L131"; Value = 6264.2724 }; @{ Cell = "N131"; Value = -16344.2724 })
    )
    "GSM" = @(
        @(@{ Cell = "H70"; Value = 0 }; @{ Cell = "I70"; Value = 0 }; @{ Cell = "K70"; Value = 0 }; @{ Cell = "M70"; Value = $null })
        @(@{ Cell = "H73"; Value = 0 }; @{ Cell = "I73"; Value = 0 }; @{ Cell = "K73"; Value = 0 }; @{ Cell = "M73"; Value = $null })
        @(@{ Cell = "H122"; Value = 0 }; @{ Cell = "I122"; Value = 0 }; @{ Cell = "K122"; Value = 0 }; @{ Cell = "M122"; Value = $null })
        @(@{ Cell = "H133"; Value = 0 }; @{ Cell = "J133"; Value = 0 }; @{ Cell = "L133"; Value = 0 }; @{ Cell = "N133"; Value = $null })
    )
    "LTW" = @(
        @(@{ Cell = "H7"; Value = 1799.6666 }; @{ Cell = "I7"; Value = 1799.6666 }; @{ Cell = "K7"; Value = 1799.6666 }; @{ Cell = "M7"; Value = -1687.6666 })
        @(@{ Cell = "H16"; Value = 636.125 }; @{ Cell = "I16"; Value = 612.7143 }; @{ Cell = "K16"; Value = 612.7143 }; @{ Cell = "M16"; Value = -442.7143 })
        @(@{ Cell = "H109"; Value = 0 }; @{ Cell = "J109"; Value = 0 }; @{ Cell = "L109"; Value = 0 }; @{ Cell = "N109"; Value = $null })
        @(@{ Cell = "H126"; Value = 1799.6666 }; @{ Cell = "I126"; Value = 1799.6666 }; @{ Cell = "K126"; Value = 5398.9998 }; @{ Cell = "M126"; Value = -2928.9998 })
    )
    "WVR" = @(
        @(@{ Cell = "H132"; Value = 2847.65 }; @{ Cell = "I132"; Value = 2283.6924 }; @{ Cell = "K132"; Value = 6851.0772 }; @{ Cell = "M132"; Value = -4321.0772 })
        @(@{ Cell = "H136"; Value = 2699.25 }; @{ Cell = "I136"; Value = 1899.5 }; @{ Cell = "K136"; Value = 5698.5 }; @{ Cell = "M136"; Value = -3148.5 })
    )
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($rowUpdates in $sheetUpdates[$sheetName]) {
        foreach ($update in $rowUpdates) {
            if ($null -eq $update.Value) {
                $ws.Range($update.Cell).ClearContents()
            } else {
                $ws.Range($update.Cell).Value = $update.Value
            }
        }
    }
}

Write-Output "Applied scheduled profit-column updates to $($sheetUpdates.Keys.Count) sheets."